$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 45, pushing the existing rows 45-77 down to 46-78.
$ws.Rows(45).Insert()

# Populate the newly inserted row 45 with the new weekly record.
$ws.Range("A45").Value = 11
$ws.Range("B45").Value = 'Vega Monumental Concepción'
$ws.Range("C45").Value = 'Bíobío'
$ws.Range("D45").Value = 44719
$ws.Range("E45").Value = 8
$ws.Range("F45").Value = 100112012
$ws.Range("G45").Value = 'Espinaca'
$ws.Range("H45").Value = 'Sin especificar'
$ws.Range("I45").Value = 'Primera'
$ws.Range("J45").Value = 100
$ws.Range("K45").Value = 7000
$ws.Range("L45").Value = 7500
$ws.Range("M45").Value = 7250
$ws.Range("N45").Value = '$/cuna 10 kilos'
$ws.Range("O45").Value = 'Región Metropolitana'
$ws.Range("P45").Value = 725
$ws.Range("Q45").Value = 10
$ws.Range("R45").Value = 'Hortaliza'
